$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 14-46 (their content is merged into rows 2-13 below)
$ws.Range("A14:A46").EntireRow.Delete()

# Rewrite rows 2-13 with the combined tuple-style strings
$ws.Range("A2").Value = "('Bear', ['Token Creature — Bear', '4/4'])"
$ws.Range("A3").Value = "('Bird', ['Token Creature — Bird', 'Flying', '3/4'])"
$ws.Range("A4").Value = "('Goblin', ['Token Creature — Goblin', '1/1'])"
$ws.Range("A5").Value = "('Morph', ['Creature', '(You can cover a face-down creature with this reminder card.', 'A card with morph can be turned face up any time for its morph cost.)', '2/2'])"
$ws.Range("A6").Value = "('Sarkhan, the Dragonspeaker Emblem', ['Emblem — Sarkhan', 'At the beginning of your draw step, draw two additional cards.', 'At the beginning of your end step, discard your hand.'])"
$ws.Range("A7").Value = "('Snake', ['Token Creature — Snake', '1/1'])"
$ws.Range("A8").Value = "('Sorin, Solemn Visitor Emblem', ['Emblem — Sorin', 'At the beginning of each opponent’s upkeep, that player sacrifices a creature.'])"
$ws.Range("A9").Value = "('Spirit', ['Token Creature — Spirit', 'Flying', '1/1'])"
$ws.Range("A10").Value = "('Spirit Warrior', ['Token Creature — Spirit Warrior', '*/*'])"
$ws.Range("A11").Value = "('Vampire', ['Token Creature — Vampire', 'Flying', '2/2'])"
$ws.Range("A12").Value = "('Warrior', ['Token Creature — Warrior', '1/1'])"
$ws.Range("A13").Value = "('Zombie', ['Token Creature — Zombie', '2/2'])"
